$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("dFBA net reactions")
$ws1.Name = "dFBA objective reactions"

$ws2 = $wb.Worksheets.Item("dFBA net species")
$ws2.Name = "dFBA objective species"
